$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 43.995596
$ws.Range("H2").Value = 131.986788
$ws.Range("I2").Value = 0.08241811124115486
$ws.Range("J2").Value = 0.08241811124115485
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 21.09934133333334
$ws.Range("N2").Value = 63.29802400000001
$ws.Range("O2").Value = 0.2917236204149438
$ws.Range("P2").Value = 0.2917236204149438
$ws.Range("Q2").Value = 928.2780971674347
$ws.Range("R2").Value = 8354.502874506912
$ws.Range("S2").Value = 0.02404330979903127
$ws.Range("T2").Value = 0.02404330979903127

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 43.995596
$ws.Range("H3").Value = 131.986788
$ws.Range("I3").Value = 0.08241811124115486
$ws.Range("J3").Value = 0.08241811124115485
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 35.81943766666667
$ws.Range("N3").Value = 107.458313
$ws.Range("O3").Value = 0.4952465516465762
$ws.Range("P3").Value = 0.4952465516465762
$ws.Range("Q3").Value = 1575.897508529849
$ws.Range("R3").Value = 14183.07757676864
$ws.Range("S3").Value = 0.04081728538540586
$ws.Range("T3").Value = 0.04081728538540586

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 43.995596
$ws.Range("H4").Value = 131.986788
$ws.Range("I4").Value = 0.08241811124115486
$ws.Range("J4").Value = 0.08241811124115485
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 15.40769666666667
$ws.Range("N4").Value = 46.22309
$ws.Range("O4").Value = 0.2130298279384801
$ws.Range("P4").Value = 0.2130298279384801
$ws.Range("Q4").Value = 677.8707978372133
$ws.Range("R4").Value = 6100.83718053492
$ws.Range("S4").Value = 0.01755751605671773
$ws.Range("T4").Value = 0.01755751605671773

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 439.8208616666666
$ws.Range("H5").Value = 1319.462585
$ws.Range("I5").Value = 0.8239280291378236
$ws.Range("J5").Value = 0.8239280291378236
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 21.09934133333334
$ws.Range("N5").Value = 63.29802400000001
$ws.Range("O5").Value = 0.2917236204149438
$ws.Range("P5").Value = 0.2917236204149438
$ws.Range("Q5").Value = 9279.930485825782
$ws.Range("R5").Value = 83519.37437243205
$ws.Range("S5").Value = 0.2403592676214352
$ws.Range("T5").Value = 0.2403592676214352

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 439.8208616666666
$ws.Range("H6").Value = 1319.462585
$ws.Range("I6").Value = 0.8239280291378236
$ws.Range("J6").Value = 0.8239280291378236
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 35.81943766666667
$ws.Range("N6").Value = 107.458313
$ws.Range("O6").Value = 0.4952465516465762
$ws.Range("P6").Value = 0.4952465516465762
$ws.Range("Q6").Value = 15754.13593896879
$ws.Range("R6").Value = 141787.2234507191
$ws.Range("S6").Value = 0.4080475152354668
$ws.Range("T6").Value = 0.4080475152354669

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 439.8208616666666
$ws.Range("H7").Value = 1319.462585
$ws.Range("I7").Value = 0.8239280291378236
$ws.Range("J7").Value = 0.8239280291378236
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 15.40769666666667
$ws.Range("N7").Value = 46.22309
$ws.Range("O7").Value = 0.2130298279384801
$ws.Range("P7").Value = 0.2130298279384801
$ws.Range("Q7").Value = 6776.62642423196
$ws.Range("R7").Value = 60989.63781808765
$ws.Range("S7").Value = 0.1755212462809216
$ws.Range("T7").Value = 0.1755212462809216

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 49.99334866666667
$ws.Range("H8").Value = 149.980046
$ws.Range("I8").Value = 0.09365385962102149
$ws.Range("J8").Value = 0.09365385962102149
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 21.09934133333334
$ws.Range("N8").Value = 63.29802400000001
$ws.Range("O8").Value = 0.2917236204149438
$ws.Range("P8").Value = 0.2917236204149438
$ws.Range("Q8").Value = 1054.826727914345
$ws.Range("R8").Value = 9493.440551229105
$ws.Range("S8").Value = 0.0273210429944773
$ws.Range("T8").Value = 0.02732104299447731

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 49.99334866666667
$ws.Range("H9").Value = 149.980046
$ws.Range("I9").Value = 0.09365385962102149
$ws.Range("J9").Value = 0.09365385962102149
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 35.81943766666667
$ws.Range("N9").Value = 107.458313
$ws.Range("O9").Value = 0.4952465516465762
$ws.Range("P9").Value = 0.4952465516465762
$ws.Range("Q9").Value = 1790.7336363136
$ws.Range("R9").Value = 16116.6027268224
$ws.Range("S9").Value = 0.04638175102570341
$ws.Range("T9").Value = 0.04638175102570342

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 49.99334866666667
$ws.Range("H10").Value = 149.980046
$ws.Range("I10").Value = 0.09365385962102149
$ws.Range("J10").Value = 0.09365385962102149
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 15.40769666666667
$ws.Range("N10").Value = 46.22309
$ws.Range("O10").Value = 0.2130298279384801
$ws.Range("P10").Value = 0.2130298279384801
$ws.Range("Q10").Value = 770.2823516069045
$ws.Range("R10").Value = 6932.541164462141
$ws.Range("S10").Value = 0.01995106560084077
$ws.Range("T10").Value = 0.01995106560084078
